# ---------------------------------------------------------------------------
# "Add files via upload" -- refresh Sheet2's edge-list / self-loop tables
# with the new graph data, and restore both sheets' scroll position back
# to the top-left of the used range (Excel had scrolled down while the
# data was being edited).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Sheet1: only the view changed (scrolled back up); selection unchanged ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("E5").Select()

# --- Sheet2: becomes the active sheet again, with new data ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

$data = New-Object 'object[,]' 68,3
$data[0,0] = 0; $data[0,1] = 1; $data[0,2] = 5
$data[1,0] = 1; $data[1,1] = 2; $data[1,2] = 11
$data[2,0] = 2; $data[2,1] = 25; $data[2,2] = 6
$data[3,0] = 25; $data[3,1] = 3; $data[3,2] = 6
$data[4,0] = 3; $data[4,1] = 26; $data[4,2] = 6
$data[5,0] = 26; $data[5,1] = 4; $data[5,2] = 6
$data[6,0] = 4; $data[6,1] = 5; $data[6,2] = 5
$data[7,0] = 5; $data[7,1] = 27; $data[7,2] = 7
$data[8,0] = 27; $data[8,1] = 6; $data[8,2] = 7
$data[9,0] = 6; $data[9,1] = 28; $data[9,2] = 7
$data[10,0] = 28; $data[10,1] = 7; $data[10,2] = 7
$data[11,0] = 7; $data[11,1] = 29; $data[11,2] = 10
$data[12,0] = 29; $data[12,1] = 8; $data[12,2] = 10
$data[13,0] = 8; $data[13,1] = 9; $data[13,2] = 5
$data[14,0] = 9; $data[14,1] = 30; $data[14,2] = 6
$data[15,0] = 30; $data[15,1] = 10; $data[15,2] = 6
$data[16,0] = 10; $data[16,1] = 31; $data[16,2] = 10
$data[17,0] = 31; $data[17,1] = 32; $data[17,2] = 10
$data[18,0] = 32; $data[18,1] = 33; $data[18,2] = 6
$data[19,0] = 33; $data[19,1] = 11; $data[19,2] = 6
$data[20,0] = 11; $data[20,1] = 12; $data[20,2] = 11
$data[21,0] = 12; $data[21,1] = 13; $data[21,2] = 5
$data[22,0] = 13; $data[22,1] = 34; $data[22,2] = 10
$data[23,0] = 34; $data[23,1] = 14; $data[23,2] = 10
$data[24,0] = 14; $data[24,1] = 35; $data[24,2] = 10
$data[25,0] = 35; $data[25,1] = 15; $data[25,2] = 10
$data[26,0] = 15; $data[26,1] = 36; $data[26,2] = 7
$data[27,0] = 36; $data[27,1] = 16; $data[27,2] = 7
$data[28,0] = 16; $data[28,1] = 37; $data[28,2] = 7
$data[29,0] = 37; $data[29,1] = 0; $data[29,2] = 7
$data[30,0] = 17; $data[30,1] = 18; $data[30,2] = 11
$data[31,0] = 18; $data[31,1] = 38; $data[31,2] = 6
$data[32,0] = 38; $data[32,1] = 19; $data[32,2] = 6
$data[33,0] = 19; $data[33,1] = 39; $data[33,2] = 6
$data[34,0] = 39; $data[34,1] = 20; $data[34,2] = 6
$data[35,0] = 20; $data[35,1] = 6; $data[35,2] = 5
$data[36,0] = 21; $data[36,1] = 22; $data[36,2] = 11
$data[37,0] = 22; $data[37,1] = 23; $data[37,2] = 6
$data[38,0] = 23; $data[38,1] = 40; $data[38,2] = 9
$data[39,0] = 40; $data[39,1] = 24; $data[39,2] = 9
$data[40,0] = 24; $data[40,1] = 7; $data[40,2] = 5
$data[41,0] = 15; $data[41,1] = 21; $data[41,2] = 5
$data[42,0] = 16; $data[42,1] = 17; $data[42,2] = 5
$data[43,0] = 0; $data[43,1] = 0; $data[43,2] = 0
$data[44,0] = 1; $data[44,1] = 1; $data[44,2] = 0
$data[45,0] = 2; $data[45,1] = 2; $data[45,2] = 0
$data[46,0] = 3; $data[46,1] = 3; $data[46,2] = 0
$data[47,0] = 4; $data[47,1] = 4; $data[47,2] = 0
$data[48,0] = 5; $data[48,1] = 5; $data[48,2] = 0
$data[49,0] = 6; $data[49,1] = 6; $data[49,2] = 0
$data[50,0] = 7; $data[50,1] = 7; $data[50,2] = 0
$data[51,0] = 8; $data[51,1] = 8; $data[51,2] = 0
$data[52,0] = 9; $data[52,1] = 9; $data[52,2] = 0
$data[53,0] = 10; $data[53,1] = 10; $data[53,2] = 0
$data[54,0] = 11; $data[54,1] = 11; $data[54,2] = 0
$data[55,0] = 12; $data[55,1] = 12; $data[55,2] = 0
$data[56,0] = 13; $data[56,1] = 13; $data[56,2] = 0
$data[57,0] = 14; $data[57,1] = 14; $data[57,2] = 0
$data[58,0] = 15; $data[58,1] = 15; $data[58,2] = 0
$data[59,0] = 16; $data[59,1] = 16; $data[59,2] = 0
$data[60,0] = 17; $data[60,1] = 17; $data[60,2] = 0
$data[61,0] = 18; $data[61,1] = 18; $data[61,2] = 0
$data[62,0] = 19; $data[62,1] = 19; $data[62,2] = 0
$data[63,0] = 20; $data[63,1] = 20; $data[63,2] = 0
$data[64,0] = 21; $data[64,1] = 21; $data[64,2] = 0
$data[65,0] = 22; $data[65,1] = 22; $data[65,2] = 0
$data[66,0] = 23; $data[66,1] = 23; $data[66,2] = 0
$data[67,0] = 24; $data[67,1] = 24; $data[67,2] = 0
$ws2.Range("A2:C69").Value = $data

# The self-loop block (Node1==Node2, Distance 0) used to sit in rows 29:53;
# it now sits 16 rows further down, in rows 45:69. Move its distinctive
# font style along with it: clear the vacated rows, and stamp the newly
# occupied rows with the same formatting.
$ws2.Range("A53:B53").Copy()
$ws2.Range("A54:B69").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("A29:B44").ClearFormats()

# Sheet2's dimensions grew, so update the recorded used range.
$ws2.Range("A44").EntireRow.Select()
